$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.864.13"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").Value = "2.123.35"
$ws.Range("E3").Value = "  +1.44%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.24"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.617"
$ws.Range("E6").Value = "  +0.69%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "62.37"
$ws.Range("E7").Value = "  +2.66%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +1.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0845"
$ws.Range("E10").Value = "  +1.06%  "
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.92"
$ws.Range("E12").Value = "  +6.14%  "
$ws.Range("D13").Value = "2.437.04"
$ws.Range("E13").Value = "  +1.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.12"
$ws.Range("E14").Value = "  +0.31%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.808"
$ws.Range("E15").Value = "  +1.39%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.51"
$ws.Range("E16").Value = "  +0.82%  "
$ws.Range("D17").Value = "2.108.63"
$ws.Range("E17").Value = "  +0.53%  "
$ws.Range("D18").Value = "38.920.13"
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.82"
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.09"
$ws.Range("E20").Value = "  +1.01%  "
$ws.Range("D21").Value = "0.0₃0847"
$ws.Range("E21").Value = "  +1.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.67"
$ws.Range("E22").Value = "  +0.59%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.41"
$ws.Range("E24").Value = "  -0.63%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.33"
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.53"
$ws.Range("E26").Value = "  +1.15%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "170.23"
$ws.Range("E27").Value = "  -0.18%  "
$ws.Range("E28").Value = "  -0.22%  "
$ws.Range("E29").Value = "  -1.78%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.48"
$ws.Range("E30").Value = "  +1.60%  "
$ws.Range("E31").Value = "  +9.40%  "
$ws.Range("E32").Value = "  +0.66%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.61"
$ws.Range("E33").Value = "  +3.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.20"
$ws.Range("E34").Value = "  +11.77%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.79"
$ws.Range("E35").Value = "  +1.85%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0615"
$ws.Range("E36").Value = "  +0.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.40"
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.54"
$ws.Range("E38").Value = "  +1.02%  "
$ws.Range("E39").Value = "  +0.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.21"
$ws.Range("E40").Value = "  -0.31%  "
$ws.Range("E41").Value = "  +3.49%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "102.61"
$ws.Range("E42").Value = "  +1.50%  "
$ws.Range("D43").Value = "1.529.59"
$ws.Range("E43").Value = "  -0.54%  "
$ws.Range("E44").Value = "  +7.28%  "
$ws.Range("E45").Value = "  -0.67%  "
$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.09"
$ws.Range("E46").Value = "  +5.88%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.79"
$ws.Range("E47").Value = "  +1.94%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0916"
$ws.Range("E48").Value = "  -0.77%  "
$ws.Range("E49").Value = "  +1.52%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.97"
$ws.Range("E50").Value = "  -0.05%  "
$ws.Range("D51").Value = "2.324.29"
$ws.Range("E51").Value = "  +1.39%  "
